# Commit: "Added in links that work by changing textContent to innerHTML"
#
# The generator's descriptions move from plain textContent to innerHTML,
# so a few description strings on Sheet2 gain an <a href="..."> link
# (and one previously-blank description is filled in). Sheet1's formulas
# recompute automatically to reflect the new text. The workbook is also
# left with Sheet2 as the active/selected sheet (cell B10 selected).

$wb  = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("Sheet2")

# Sacred Space -> link to sacredspace.ie
$ws2.Range("B2").Value = "Go to <a href='sacredspace.ie'>Sacred Space</a>."

# Examen -> previously empty, now links to the Ignatian Examen article
$ws2.Range("B4").Value = "Go to <a href='https://www.jesuits.org/spirituality/the-ignatian-examen/'>Ignatian Examen</a>."

# Imaginative Prayer -> expanded description text
$ws2.Range("B9").Value = "Imagine one of the visible scenes in Scripture."

# Match the saved view state: Sheet2 active, B10 selected
$ws2.Activate()
$ws2.Range("B10").Select()
